$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '71.768.39'
$ws.Cells.Item(2, 5).Value = '  +4.84%  '
$ws.Cells.Item(3, 4).Value = '4.047.44'
$ws.Cells.Item(3, 5).Value = '  +5.01%  '
$ws.Cells.Item(4, 5).Value = '  +0.05%  '
$cell = $ws.Cells.Item(5, 4)
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '534.58'
$cell.Style = $origStyle
$ws.Cells.Item(5, 5).Value = '  +2.65%  '
$cell = $ws.Cells.Item(6, 4)
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '152.97'
$cell.Style = $origStyle
$ws.Cells.Item(6, 5).Value = '  +8.97%  '
$ws.Cells.Item(7, 5).Value = '  +14.30%  '
$ws.Cells.Item(8, 5).Value = '  +0.03%  '
$cell = $ws.Cells.Item(9, 4)
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.760'
$cell.Style = $origStyle
$ws.Cells.Item(9, 5).Value = '  +7.03%  '
$cell = $ws.Cells.Item(10, 4)
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.175'
$cell.Style = $origStyle
$ws.Cells.Item(10, 5).Value = '  +5.77%  '
$cell = $ws.Cells.Item(11, 4)
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.0000332'
$cell.Style = $origStyle
$ws.Cells.Item(11, 5).Value = '  +4.82%  '
$cell = $ws.Cells.Item(12, 4)
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '48.88'
$cell.Style = $origStyle
$ws.Cells.Item(12, 5).Value = '  +18.20%  '
$cell = $ws.Cells.Item(13, 4)
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '10.96'
$cell.Style = $origStyle
$ws.Cells.Item(13, 5).Value = '  +6.48%  '
$ws.Cells.Item(14, 4).Value = '4.694.03'
$ws.Cells.Item(14, 5).Value = '  +5.09%  '
$ws.Cells.Item(15, 4).Value = '4.041.84'
$ws.Cells.Item(15, 5).Value = '  +5.00%  '
$cell = $ws.Cells.Item(16, 4)
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '14.42'
$cell.Style = $origStyle
$ws.Cells.Item(16, 5).Value = '  +2.44%  '
$cell = $ws.Cells.Item(17, 4)
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '21.05'
$cell.Style = $origStyle
$ws.Cells.Item(17, 5).Value = '  -0.89%  '
$ws.Cells.Item(18, 5).Value = '  +2.35%  '
$ws.Cells.Item(19, 5).Value = '  +0.20%  '
$ws.Cells.Item(20, 4).Value = '71.788.14'
$ws.Cells.Item(20, 5).Value = '  +4.86%  '
$cell = $ws.Cells.Item(21, 4)
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '435.75'
$cell.Style = $origStyle
$ws.Cells.Item(21, 5).Value = '  +4.95%  '
$ws.Cells.Item(22, 5).Value = '  +7.83%  '
$cell = $ws.Cells.Item(23, 4)
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '99.57'
$cell.Style = $origStyle
$ws.Cells.Item(23, 5).Value = '  +15.03%  '
$cell = $ws.Cells.Item(24, 4)
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '14.82'
$cell.Style = $origStyle
$ws.Cells.Item(24, 5).Value = '  +6.19%  '
$cell = $ws.Cells.Item(25, 4)
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '4.22'
$cell.Style = $origStyle
$ws.Cells.Item(25, 5).Value = '  +6.79%  '
$cell = $ws.Cells.Item(26, 4)
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '11.42'
$cell.Style = $origStyle
$ws.Cells.Item(26, 5).Value = '  -0.80%  '
$cell = $ws.Cells.Item(27, 4)
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '10.94'
$cell.Style = $origStyle
$ws.Cells.Item(27, 5).Value = '  +4.13%  '
$cell = $ws.Cells.Item(28, 4)
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '37.28'
$cell.Style = $origStyle
$ws.Cells.Item(28, 5).Value = '  +5.39%  '
$cell = $ws.Cells.Item(29, 4)
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '5.82'
$cell.Style = $origStyle
$ws.Cells.Item(29, 5).Value = '  +2.97%  '
$cell = $ws.Cells.Item(30, 4)
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '3.54'
$cell.Style = $origStyle
$ws.Cells.Item(30, 5).Value = '  +27.54%  '
$cell = $ws.Cells.Item(31, 4)
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '13.71'
$cell.Style = $origStyle
$ws.Cells.Item(31, 5).Value = '  +4.16%  '
$ws.Cells.Item(32, 5).Value = '  +6.48%  '
$cell = $ws.Cells.Item(33, 4)
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '680.38'
$cell.Style = $origStyle
$ws.Cells.Item(33, 5).Value = '  +0.39%  '
$cell = $ws.Cells.Item(34, 4)
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '6.80'
$cell.Style = $origStyle
$ws.Cells.Item(34, 5).Value = '  +2.67%  '
$cell = $ws.Cells.Item(35, 4)
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '66.74'
$cell.Style = $origStyle
$ws.Cells.Item(35, 5).Value = '  +1.09%  '
$cell = $ws.Cells.Item(36, 4)
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '43.10'
$cell.Style = $origStyle
$ws.Cells.Item(36, 5).Value = '  +9.04%  '
$ws.Cells.Item(37, 5).Value = '  -4.58%  '
$ws.Cells.Item(38, 5).Value = '  +6.59%  '
$ws.Cells.Item(39, 4).Value = '0.0₃0860'
$ws.Cells.Item(39, 5).Value = '  +4.18%  '
$ws.Cells.Item(40, 5).Value = '  -1.97%  '
$ws.Cells.Item(41, 5).Value = '  +0.06%  '
$ws.Cells.Item(42, 5).Value = '  +5.57%  '
$cell = $ws.Cells.Item(43, 4)
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.999'
$cell.Style = $origStyle
$ws.Cells.Item(43, 5).Value = '  -0.09%  '
$cell = $ws.Cells.Item(44, 4)
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '3.18'
$cell.Style = $origStyle
$ws.Cells.Item(44, 5).Value = '  +2.50%  '
$cell = $ws.Cells.Item(45, 4)
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.153'
$cell.Style = $origStyle
$ws.Cells.Item(45, 5).Value = '  +9.35%  '
$cell = $ws.Cells.Item(46, 4)
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.73'
$cell.Style = $origStyle
$ws.Cells.Item(46, 5).Value = '  -2.20%  '
$cell = $ws.Cells.Item(47, 4)
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '3.40'
$cell.Style = $origStyle
$ws.Cells.Item(47, 5).Value = '  +0.32%  '
$ws.Cells.Item(48, 5).Value = '  +12.54%  '
$ws.Cells.Item(49, 5).Value = '  +2.07%  '
$cell = $ws.Cells.Item(50, 4)
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '3.39'
$cell.Style = $origStyle
$ws.Cells.Item(50, 5).Value = '  +4.37%  '
$cell = $ws.Cells.Item(51, 4)
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.000276'
$cell.Style = $origStyle
$ws.Cells.Item(51, 5).Value = '  +3.80%  '
